$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 44796
$ws.Range("B15").Value = "Added PTC-fuse and power-on indication-LED."

$ws.Range("B16").Select()
